$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Orders")
$summary = $wb.Worksheets.Item("Summary")

# New order rows (A1:L11 -> A1:L21). Columns A/F hold numbers written by the
# source system as text, so force a Text number format before assigning the
# value (otherwise COM auto-coerces the numeric-looking string to a Number).
function Set-TextValue($range, [string]$value) {
    $range.NumberFormat = "@"
    $range.Value = $value
}

Set-TextValue $ws.Range("A12") "2"
$ws.Range("C12").Value = "327_文竹_asparagus fern_undefined_1bunch"
Set-TextValue $ws.Range("F12") "10"

$ws.Range("C13").Value = "522_山归来绿_Smilax china_undefined_1bunch"
Set-TextValue $ws.Range("F13") "5"

$ws.Range("C14").Value = "328_卢荀草_undefined_undefined_1bunch"
Set-TextValue $ws.Range("F14") "15"

$ws.Range("C15").Value = "328_卢荀草_undefined_undefined_1bunch"
Set-TextValue $ws.Range("F15") "10"

$ws.Range("C16").Value = "100_绣球单瓣白_Hydrangea White S_Hydrangea L._1stem"
Set-TextValue $ws.Range("F16") "30"

$ws.Range("C17").Value = "100_绣球单瓣白_Hydrangea White S_Hydrangea L._1stem"
Set-TextValue $ws.Range("F17") "10"

$ws.Range("C18").Value = "106_绣球单瓣粉_Hydrangea Pink S_Hydrangea L._1stem"
Set-TextValue $ws.Range("F18") "30"

$ws.Range("C19").Value = "816_山里红_undefined_undefined_1bunch"
Set-TextValue $ws.Range("F19") "15"

$ws.Range("C20").Value = "479_绿灵草_lepidium_undefined_1bunch"
Set-TextValue $ws.Range("F20") "15"

Set-TextValue $ws.Range("A21") "3"

# Summary!G2 gains more digits appended to the packed code; keep it text so
# the leading zero and full digit string survive.
Set-TextValue $summary.Range("G2") "0201020105555530105151030103015150"
